$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each new expense row: Category, Description, Date, Amount, (optional) Account.
# These mirror the app's in-memory ArrayList of expenses per account, which now
# also gets appended to the bank-statement worksheet whenever a new expense is
# added (so the pie chart that reads this sheet stays in sync).
$newRows = @(
    @("Food",           "sdfg",      "2023-03-08", "10000.0", $null),
    @("Food",           "dgrdgf",    "2023-03-11", "1500.0",  $null),
    @("Food",           "adms",      "2023-03-11", "150.0",   $null),
    @("Transportation", "dfms",      "2023-03-11", "150.0",   $null),
    @("Food",           "asws",      "2023-03-11", "150.0",   $null),
    @("Entertainment",  "KSLKFVASZ", "2023-03-11", "99.0",    "Savings"),
    @("Other",          "gthfthfg",  "2023-03-11", "4999.0",  "Savings"),
    @("Rent",           "gtfg",      "2023-03-11", "545.0",   "Savings")
)

$startRow = 25
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    # Date / amount columns look numeric to the app, so they'd otherwise be
    # auto-converted to a date serial / number. Force text while typing the
    # value in, then drop the temporary number format again so the cell is
    # left with plain literal text, same as the rest of the sheet.
    $cellC = $ws.Cells.Item($r, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $row[2]
    $cellC.ClearFormats()

    $cellD = $ws.Cells.Item($r, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $row[3]
    $cellD.ClearFormats()

    if ($row[4] -ne $null) {
        $ws.Cells.Item($r, 5).Value = $row[4]
    }
}
